$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.036.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.41%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.82"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.38%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7021"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.96"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.12%  "

# Row 7
$ws.Range("E7").Value = "  +0.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3034"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.35%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07461"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.85%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.60%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08124"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.88%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.51"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.60%  "

# Row 13
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7249"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.65%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.231"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.92"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.18%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.157.58"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.781"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.74"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.42%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007660"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.02"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.31%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.105.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.31%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.27%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.547"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.85%  "

# Row 25
$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1459"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.79%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.938"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.95"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.378"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.511"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.07%  "

# Row 32
$ws.Range("E32").Value = "  -2.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.998"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05145"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.183"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.84%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7053"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.49%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.007"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.644"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01856"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.99%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.666"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.69%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9013"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.976"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4278"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.061.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.31%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").ClearFormats()

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.747"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.75%  "

# Row 49
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.986.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.046"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.79%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.148"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.64%  "
